# Update the "Price" (D) and "Volume(1h)" (E) columns with refreshed
# quotes, matching the latest symbol-list run (GitHub Actions refresh).
#
# Source cells are stored as literal text (e.g. "328.25", "-1.02%") rather
# than numbers, so each target cell is first switched to the Text number
# format ("@") before the new literal is written — this stops Excel's
# normal type-coercion from turning "328.48" into the number 328.48 (and
# losing things like trailing zeros) or turning "-0.94%" into a numeric
# percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# row -> @(new Price, new Volume(1h))  ("" means Price is unchanged)
$updates = [ordered]@{
    2  = @("328.48",         "-0.94%")
    3  = @("44.05",          "6.69%")
    4  = @("5.547",          "-2.78%")
    5  = @("0.08093",        "-3.79%")
    6  = @("8.675",          "-1.44%")
    7  = @("1.919",          "-3.20%")
    8  = @("4.301",          "-4.78%")
    10 = @("0.9468",         "2.34%")
    11 = @("0.1190",         "-4.22%")
    12 = @("0.1896",         "-3.15%")
    13 = @("0.09659",        "3.22%")
    14 = @("0.04072",        "2.46%")
    15 = @("0.1066",         "0.22%")
    16 = @("0.001269",       "-3.03%")
    17 = @("0.005957",       "-2.56%")
    18 = @("3.580",          "4.22%")
    20 = @("8.515",          "-6.50%")
    21 = @("0.1363",         "-0.91%")
    22 = @("0.2589",         "-1.62%")
    23 = @("0.04352",        "-1.70%")
    24 = @("0.001238",       "-0.59%")
    25 = @("0.004634",       "6.04%")
    26 = @("0.0001231",      "2.83%")
    27 = @("0.0003994",      "-0.34%")
    39 = @("0.02670",        "-4.58%")
    40 = @("0.05621",        "1.78%")
    41 = @("",               "26.04%")
    42 = @("0.007639",       "-3.48%")
    43 = @("0.1395",         "-2.49%")
    44 = @("0.002110",       "-3.53%")
    45 = @("0.009434",       "-9.21%")
    46 = @("0.00007083",     "-1.46%")
    47 = @("0.00000000751",  "-0.47%")
    48 = @("0.003445",       "-0.72%")
    49 = @("0.002273",       "-0.67%")
    50 = @("0.00002103",     "-0.47%")
    51 = @("0.0002002",      "-0.47%")
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $price = $pair[0]
    $volume = $pair[1]

    if ($price -ne "") {
        Set-TextValue $ws.Range("D$row") $price
    }
    Set-TextValue $ws.Range("E$row") $volume
}
